# Update cryptos list (price + volume-change figures, and row 43/44 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.983.42"
$ws.Range("E2").Value = "  -7.17%  "
$ws.Range("D3").Value = "'3.473.40"
$ws.Range("E3").Value = "  -3.43%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'389.34"
$ws.Range("E5").Value = "  -6.94%  "
$ws.Range("D6").Value = "'120.29"
$ws.Range("E6").Value = "  -8.04%  "
$ws.Range("D7").Value = "'3.467.88"
$ws.Range("E7").Value = "  -3.41%  "
$ws.Range("E8").Value = "  -10.16%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'0.671"
$ws.Range("E10").Value = "  -12.89%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -16.29%  "
$ws.Range("D12").Value = "'0.0000329"
$ws.Range("E12").Value = "  -3.80%  "
$ws.Range("D13").Value = "'38.46"
$ws.Range("E13").Value = "  -9.98%  "
$ws.Range("D14").Value = "'4.040.81"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "'9.12"
$ws.Range("E15").Value = "  -8.79%  "
$ws.Range("D16").Value = "'0.136"
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "'3.496.15"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "'18.60"
$ws.Range("E18").Value = "  -9.25%  "
$ws.Range("D19").Value = "'12.44"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").Value = "'63.064.41"
$ws.Range("E20").Value = "  -6.81%  "
$ws.Range("E21").Value = "  -11.95%  "
$ws.Range("D22").Value = "'391.51"
$ws.Range("E22").Value = "  -15.52%  "
$ws.Range("D23").Value = "'13.82"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "'80.53"
$ws.Range("E24").Value = "  -9.18%  "
$ws.Range("D25").Value = "'2.85"
$ws.Range("E25").Value = "  -9.43%  "
$ws.Range("D26").Value = "'33.22"
$ws.Range("E26").Value = "  -7.04%  "
$ws.Range("D27").Value = "'5.15"
$ws.Range("E27").Value = "  +6.08%  "
$ws.Range("D28").Value = "'2.97"
$ws.Range("E28").Value = "  -11.91%  "
$ws.Range("D29").Value = "'8.72"
$ws.Range("E29").Value = "  -15.36%  "
$ws.Range("D30").Value = "'11.73"
$ws.Range("E30").Value = "  -5.89%  "
$ws.Range("D31").Value = "'2.59"
$ws.Range("E31").Value = "  -7.14%  "
$ws.Range("E32").Value = "  -7.52%  "
$ws.Range("D33").Value = "'6.71"
$ws.Range("E33").Value = "  -10.27%  "
$ws.Range("D34").Value = "'0.148"
$ws.Range("E34").Value = "  -9.21%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "'36.35"
$ws.Range("E36").Value = "  -12.62%  "
$ws.Range("D37").Value = "'53.06"
$ws.Range("E37").Value = "  -6.67%  "
$ws.Range("D38").Value = "'0.0433"
$ws.Range("E38").Value = "  -12.82%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "'2.69"
$ws.Range("E40").Value = "  +14.72%  "
$ws.Range("D41").Value = "'0.0₃0633"
$ws.Range("E41").Value = "  -11.75%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'140.91"
$ws.Range("E43").Value = "  -5.08%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'3.05"
$ws.Range("E44").Value = "  +12.35%  "
$ws.Range("D45").Value = "'2.73"
$ws.Range("E45").Value = "  -10.32%  "
$ws.Range("E46").Value = "  -7.14%  "
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").Value = "'24.42"
$ws.Range("E48").Value = "  +12.38%  "
$ws.Range("D49").Value = "'2.44"
$ws.Range("E49").Value = "  -10.73%  "
$ws.Range("D50").Value = "'3.97"
$ws.Range("E50").Value = "  -8.31%  "
$ws.Range("D51").Value = "'0.273"
$ws.Range("E51").Value = "  -12.68%  "
